# Reload the exam/grade data: re-sorted (alphabetically by last name),
# names normalized to lower-case, and the exact/rounded grade columns
# (C: "Note Exakt", D: "Note Gerundet") populated from the freshly
# loaded exam results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("asdf",     "marlene", 5.175109999999999, 5.25),
    @("kohler",   "alina",   1,                 1),
    @("kohler",   "nina",    5.37826,            5.5),
    @("matumona", "noe",     6,                  6),
    @("matumona", "nina",    5.82464,            5.75),
    @("sarman",   "dominik", 4.78031,            4.75),
    @("zillig",   "nicolas", 3.94643,            4)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
